# Auto-generated-intent PowerShell script (hand-verified) to apply the logging-config upgrade diff
$wb = $excel.ActiveWorkbook

# --- Update existing row 181 timestamp (A column) on sheets 2-4 ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(181,1).Value = 45967.4921412037

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(181,1).Value = 45967.4921412037

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(181,1).Value = 45967.4921412037

# --- Sheet 1: append rows 182-189 ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(182,1).Value = 45968.49288194445
$ws.Cells.Item(182,2).Value = "0x01,0x7c"
$ws.Cells.Item(182,3).Value = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
$ws.Cells.Item(182,4).Value = "0x00,0xA1"
$ws.Cells.Item(182,5).Value = "0xf"
$ws.Cells.Item(182,6).Value = 380
$ws.Cells.Item(182,7).Value = [double]"7.598631275147109e+23"
$ws.Cells.Item(182,8).Value = 180
$ws.Cells.Item(182,9).Value = 15
$ws.Cells.Item(182,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()
$ws.Cells.Item(183,1).Value = 45969.49362268519
$ws.Cells.Item(183,2).Value = "0x01,0x7c"
$ws.Cells.Item(183,3).Value = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
$ws.Cells.Item(183,4).Value = "0x00,0x10C"
$ws.Cells.Item(183,5).Value = "0xf"
$ws.Cells.Item(183,6).Value = 380
$ws.Cells.Item(183,7).Value = [double]"7.598631275147109e+23"
$ws.Cells.Item(183,8).Value = 180
$ws.Cells.Item(183,9).Value = 15
$ws.Cells.Item(183,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()
$ws.Cells.Item(184,1).Value = 45970.49436342593
$ws.Cells.Item(184,2).Value = "0x01,0x7c"
$ws.Cells.Item(184,3).Value = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
$ws.Cells.Item(184,4).Value = "0x00,0xA2"
$ws.Cells.Item(184,5).Value = "0xf"
$ws.Cells.Item(184,6).Value = 380
$ws.Cells.Item(184,7).Value = [double]"7.598631275147109e+23"
$ws.Cells.Item(184,8).Value = 180
$ws.Cells.Item(184,9).Value = 15
$ws.Cells.Item(184,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()
$ws.Cells.Item(185,1).Value = 45971.49510416666
$ws.Cells.Item(185,2).Value = "0x01,0x7c"
$ws.Cells.Item(185,3).Value = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
$ws.Cells.Item(185,4).Value = "0x00,0x11C"
$ws.Cells.Item(185,5).Value = "0xf"
$ws.Cells.Item(185,6).Value = 380
$ws.Cells.Item(185,7).Value = [double]"7.598631275147109e+23"
$ws.Cells.Item(185,8).Value = 180
$ws.Cells.Item(185,9).Value = 15
$ws.Cells.Item(185,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()
$ws.Cells.Item(186,1).Value = 45972.4958449074
$ws.Cells.Item(186,2).Value = "0x01,0x7c"
$ws.Cells.Item(186,3).Value = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
$ws.Cells.Item(186,4).Value = "0x00,0xA3"
$ws.Cells.Item(186,5).Value = "0xf"
$ws.Cells.Item(186,6).Value = 380
$ws.Cells.Item(186,7).Value = [double]"7.598631275147109e+23"
$ws.Cells.Item(186,8).Value = 180
$ws.Cells.Item(186,9).Value = 15
$ws.Cells.Item(186,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()
$ws.Cells.Item(187,1).Value = 45973.49658564815
$ws.Cells.Item(187,2).Value = "0x01,0x7c"
$ws.Cells.Item(187,3).Value = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
$ws.Cells.Item(187,4).Value = "0x00,0x12C"
$ws.Cells.Item(187,5).Value = "0xf"
$ws.Cells.Item(187,6).Value = 380
$ws.Cells.Item(187,7).Value = [double]"7.598631275147109e+23"
$ws.Cells.Item(187,8).Value = 180
$ws.Cells.Item(187,9).Value = 15
$ws.Cells.Item(187,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()
$ws.Cells.Item(188,1).Value = 45974.49732638889
$ws.Cells.Item(188,2).Value = "0x01,0x7c"
$ws.Cells.Item(188,3).Value = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
$ws.Cells.Item(188,4).Value = "0x00,0xA4"
$ws.Cells.Item(188,5).Value = "0xf"
$ws.Cells.Item(188,6).Value = 380
$ws.Cells.Item(188,7).Value = [double]"7.598631275147109e+23"
$ws.Cells.Item(188,8).Value = 180
$ws.Cells.Item(188,9).Value = 15
$ws.Cells.Item(188,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()
$ws.Cells.Item(189,1).Value = 45975.49806712963
$ws.Cells.Item(189,2).Value = "0x01,0x7c"
$ws.Cells.Item(189,3).Value = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
$ws.Cells.Item(189,4).Value = "0x00,0x13C"
$ws.Cells.Item(189,5).Value = "0xf"
$ws.Cells.Item(189,6).Value = 380
$ws.Cells.Item(189,7).Value = [double]"7.598631275147109e+23"
$ws.Cells.Item(189,8).Value = 180
$ws.Cells.Item(189,9).Value = 15
$ws.Cells.Item(189,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()

# --- Sheet 2: append rows 182-189 ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(182,1).Value = 45968.49288194445
$ws.Cells.Item(182,2).Value = "0x01,0x90"
$ws.Cells.Item(182,3).Value = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
$ws.Cells.Item(182,4).Value = "0x00,0x10C"
$ws.Cells.Item(182,5).Value = "0xe"
$ws.Cells.Item(182,6).Value = 400
$ws.Cells.Item(182,7).Value = [double]"5.68432987514711e+23"
$ws.Cells.Item(182,8).Value = 164
$ws.Cells.Item(182,9).Value = 14
$ws.Cells.Item(182,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()
$ws.Cells.Item(183,1).Value = 45969.49362268519
$ws.Cells.Item(183,2).Value = "0x01,0x90"
$ws.Cells.Item(183,3).Value = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
$ws.Cells.Item(183,4).Value = "0x00,0x99"
$ws.Cells.Item(183,5).Value = "0xe"
$ws.Cells.Item(183,6).Value = 400
$ws.Cells.Item(183,7).Value = [double]"5.68432987514711e+23"
$ws.Cells.Item(183,8).Value = 164
$ws.Cells.Item(183,9).Value = 14
$ws.Cells.Item(183,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()
$ws.Cells.Item(184,1).Value = 45970.49436342593
$ws.Cells.Item(184,2).Value = "0x01,0x90"
$ws.Cells.Item(184,3).Value = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
$ws.Cells.Item(184,4).Value = "0x00,0x11C"
$ws.Cells.Item(184,5).Value = "0xe"
$ws.Cells.Item(184,6).Value = 400
$ws.Cells.Item(184,7).Value = [double]"5.68432987514711e+23"
$ws.Cells.Item(184,8).Value = 164
$ws.Cells.Item(184,9).Value = 14
$ws.Cells.Item(184,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()
$ws.Cells.Item(185,1).Value = 45971.49510416666
$ws.Cells.Item(185,2).Value = "0x01,0x90"
$ws.Cells.Item(185,3).Value = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
$ws.Cells.Item(185,4).Value = "0x00,0x100"
$ws.Cells.Item(185,5).Value = "0xe"
$ws.Cells.Item(185,6).Value = 400
$ws.Cells.Item(185,7).Value = [double]"5.68432987514711e+23"
$ws.Cells.Item(185,8).Value = 164
$ws.Cells.Item(185,9).Value = 14
$ws.Cells.Item(185,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()
$ws.Cells.Item(186,1).Value = 45972.4958449074
$ws.Cells.Item(186,2).Value = "0x01,0x90"
$ws.Cells.Item(186,3).Value = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
$ws.Cells.Item(186,4).Value = "0x00,0x12C"
$ws.Cells.Item(186,5).Value = "0xe"
$ws.Cells.Item(186,6).Value = 400
$ws.Cells.Item(186,7).Value = [double]"5.68432987514711e+23"
$ws.Cells.Item(186,8).Value = 164
$ws.Cells.Item(186,9).Value = 14
$ws.Cells.Item(186,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()
$ws.Cells.Item(187,1).Value = 45973.49658564815
$ws.Cells.Item(187,2).Value = "0x01,0x90"
$ws.Cells.Item(187,3).Value = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
$ws.Cells.Item(187,4).Value = "0x00,0x101"
$ws.Cells.Item(187,5).Value = "0xe"
$ws.Cells.Item(187,6).Value = 400
$ws.Cells.Item(187,7).Value = [double]"5.68432987514711e+23"
$ws.Cells.Item(187,8).Value = 164
$ws.Cells.Item(187,9).Value = 14
$ws.Cells.Item(187,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()
$ws.Cells.Item(188,1).Value = 45974.49732638889
$ws.Cells.Item(188,2).Value = "0x01,0x90"
$ws.Cells.Item(188,3).Value = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
$ws.Cells.Item(188,4).Value = "0x00,0x13C"
$ws.Cells.Item(188,5).Value = "0xe"
$ws.Cells.Item(188,6).Value = 400
$ws.Cells.Item(188,7).Value = [double]"5.68432987514711e+23"
$ws.Cells.Item(188,8).Value = 164
$ws.Cells.Item(188,9).Value = 14
$ws.Cells.Item(188,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()
$ws.Cells.Item(189,1).Value = 45975.49806712963
$ws.Cells.Item(189,2).Value = "0x01,0x90"
$ws.Cells.Item(189,3).Value = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
$ws.Cells.Item(189,4).Value = "0x00,0x102"
$ws.Cells.Item(189,5).Value = "0xe"
$ws.Cells.Item(189,6).Value = 400
$ws.Cells.Item(189,7).Value = [double]"5.68432987514711e+23"
$ws.Cells.Item(189,8).Value = 164
$ws.Cells.Item(189,9).Value = 14
$ws.Cells.Item(189,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()

# --- Sheet 3: append rows 182-189 ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(182,1).Value = 45968.49288194445
$ws.Cells.Item(182,2).Value = "0x00,0x6e"
$ws.Cells.Item(182,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item(182,4).Value = "0x00,0x5C"
$ws.Cells.Item(182,5).Value = "0x3"
$ws.Cells.Item(182,6).Value = 110
$ws.Cells.Item(182,7).Value = [double]"5.68631262647114e+23"
$ws.Cells.Item(182,8).Value = 75
$ws.Cells.Item(182,9).Value = 3
$ws.Cells.Item(182,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()
$ws.Cells.Item(183,1).Value = 45969.49362268519
$ws.Cells.Item(183,2).Value = "0x00,0x6e"
$ws.Cells.Item(183,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item(183,4).Value = "0x00,0x5B"
$ws.Cells.Item(183,5).Value = "0x3"
$ws.Cells.Item(183,6).Value = 110
$ws.Cells.Item(183,7).Value = [double]"5.68631262647114e+23"
$ws.Cells.Item(183,8).Value = 75
$ws.Cells.Item(183,9).Value = 3
$ws.Cells.Item(183,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()
$ws.Cells.Item(184,1).Value = 45970.49436342593
$ws.Cells.Item(184,2).Value = "0x00,0x6e"
$ws.Cells.Item(184,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item(184,4).Value = "0x00,0x6C"
$ws.Cells.Item(184,5).Value = "0x3"
$ws.Cells.Item(184,6).Value = 110
$ws.Cells.Item(184,7).Value = [double]"5.68631262647114e+23"
$ws.Cells.Item(184,8).Value = 75
$ws.Cells.Item(184,9).Value = 3
$ws.Cells.Item(184,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()
$ws.Cells.Item(185,1).Value = 45971.49510416666
$ws.Cells.Item(185,2).Value = "0x00,0x6e"
$ws.Cells.Item(185,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item(185,4).Value = "0x00,0x6B"
$ws.Cells.Item(185,5).Value = "0x3"
$ws.Cells.Item(185,6).Value = 110
$ws.Cells.Item(185,7).Value = [double]"5.68631262647114e+23"
$ws.Cells.Item(185,8).Value = 75
$ws.Cells.Item(185,9).Value = 3
$ws.Cells.Item(185,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()
$ws.Cells.Item(186,1).Value = 45972.4958449074
$ws.Cells.Item(186,2).Value = "0x00,0x6e"
$ws.Cells.Item(186,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item(186,4).Value = "0x00,0x7C"
$ws.Cells.Item(186,5).Value = "0x3"
$ws.Cells.Item(186,6).Value = 110
$ws.Cells.Item(186,7).Value = [double]"5.68631262647114e+23"
$ws.Cells.Item(186,8).Value = 75
$ws.Cells.Item(186,9).Value = 3
$ws.Cells.Item(186,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()
$ws.Cells.Item(187,1).Value = 45973.49658564815
$ws.Cells.Item(187,2).Value = "0x00,0x6e"
$ws.Cells.Item(187,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item(187,4).Value = "0x00,0x7B"
$ws.Cells.Item(187,5).Value = "0x3"
$ws.Cells.Item(187,6).Value = 110
$ws.Cells.Item(187,7).Value = [double]"5.68631262647114e+23"
$ws.Cells.Item(187,8).Value = 75
$ws.Cells.Item(187,9).Value = 3
$ws.Cells.Item(187,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()
$ws.Cells.Item(188,1).Value = 45974.49732638889
$ws.Cells.Item(188,2).Value = "0x00,0x6e"
$ws.Cells.Item(188,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item(188,4).Value = "0x00,0x8C"
$ws.Cells.Item(188,5).Value = "0x3"
$ws.Cells.Item(188,6).Value = 110
$ws.Cells.Item(188,7).Value = [double]"5.68631262647114e+23"
$ws.Cells.Item(188,8).Value = 75
$ws.Cells.Item(188,9).Value = 3
$ws.Cells.Item(188,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()
$ws.Cells.Item(189,1).Value = 45975.49806712963
$ws.Cells.Item(189,2).Value = "0x00,0x6e"
$ws.Cells.Item(189,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item(189,4).Value = "0x00,0x8B"
$ws.Cells.Item(189,5).Value = "0x3"
$ws.Cells.Item(189,6).Value = 110
$ws.Cells.Item(189,7).Value = [double]"5.68631262647114e+23"
$ws.Cells.Item(189,8).Value = 75
$ws.Cells.Item(189,9).Value = 3
$ws.Cells.Item(189,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()

# --- Sheet 4: append rows 182-189 ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(182,1).Value = 45968.49288194445
$ws.Cells.Item(182,2).Value = "0x00,0x6e"
$ws.Cells.Item(182,3).Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws.Cells.Item(182,4).Value = "0x00,0x45"
$ws.Cells.Item(182,5).Value = "0x3"
$ws.Cells.Item(182,6).Value = 110
$ws.Cells.Item(182,7).Value = [double]"9.85046333984776e+23"
$ws.Cells.Item(182,8).Value = 70
$ws.Cells.Item(182,9).Value = 3
$ws.Cells.Item(182,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()
$ws.Cells.Item(183,1).Value = 45969.49362268519
$ws.Cells.Item(183,2).Value = "0x00,0x6e"
$ws.Cells.Item(183,3).Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws.Cells.Item(183,4).Value = "0x00,0x44"
$ws.Cells.Item(183,5).Value = "0x3"
$ws.Cells.Item(183,6).Value = 110
$ws.Cells.Item(183,7).Value = [double]"9.85046333984776e+23"
$ws.Cells.Item(183,8).Value = 70
$ws.Cells.Item(183,9).Value = 3
$ws.Cells.Item(183,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()
$ws.Cells.Item(184,1).Value = 45970.49436342593
$ws.Cells.Item(184,2).Value = "0x00,0x6e"
$ws.Cells.Item(184,3).Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws.Cells.Item(184,4).Value = "0x00,0x43"
$ws.Cells.Item(184,5).Value = "0x3"
$ws.Cells.Item(184,6).Value = 110
$ws.Cells.Item(184,7).Value = [double]"9.85046333984776e+23"
$ws.Cells.Item(184,8).Value = 70
$ws.Cells.Item(184,9).Value = 3
$ws.Cells.Item(184,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()
$ws.Cells.Item(185,1).Value = 45971.49510416666
$ws.Cells.Item(185,2).Value = "0x00,0x6e"
$ws.Cells.Item(185,3).Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws.Cells.Item(185,4).Value = "0x00,0x42"
$ws.Cells.Item(185,5).Value = "0x3"
$ws.Cells.Item(185,6).Value = 110
$ws.Cells.Item(185,7).Value = [double]"9.85046333984776e+23"
$ws.Cells.Item(185,8).Value = 70
$ws.Cells.Item(185,9).Value = 3
$ws.Cells.Item(185,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()
$ws.Cells.Item(186,1).Value = 45972.4958449074
$ws.Cells.Item(186,2).Value = "0x00,0x6e"
$ws.Cells.Item(186,3).Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws.Cells.Item(186,4).Value = "0x00,0x41"
$ws.Cells.Item(186,5).Value = "0x3"
$ws.Cells.Item(186,6).Value = 110
$ws.Cells.Item(186,7).Value = [double]"9.85046333984776e+23"
$ws.Cells.Item(186,8).Value = 70
$ws.Cells.Item(186,9).Value = 3
$ws.Cells.Item(186,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()
$ws.Cells.Item(187,1).Value = 45973.49658564815
$ws.Cells.Item(187,2).Value = "0x00,0x6e"
$ws.Cells.Item(187,3).Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws.Cells.Item(187,4).Value = "0x00,0x40"
$ws.Cells.Item(187,5).Value = "0x3"
$ws.Cells.Item(187,6).Value = 110
$ws.Cells.Item(187,7).Value = [double]"9.85046333984776e+23"
$ws.Cells.Item(187,8).Value = 70
$ws.Cells.Item(187,9).Value = 3
$ws.Cells.Item(187,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()
$ws.Cells.Item(188,1).Value = 45974.49732638889
$ws.Cells.Item(188,2).Value = "0x00,0x6e"
$ws.Cells.Item(188,3).Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws.Cells.Item(188,4).Value = "0x00,0x39"
$ws.Cells.Item(188,5).Value = "0x3"
$ws.Cells.Item(188,6).Value = 110
$ws.Cells.Item(188,7).Value = [double]"9.85046333984776e+23"
$ws.Cells.Item(188,8).Value = 69
$ws.Cells.Item(188,9).Value = 3
$ws.Cells.Item(188,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()
$ws.Cells.Item(189,1).Value = 45975.49806712963
$ws.Cells.Item(189,2).Value = "0x00,0x6e"
$ws.Cells.Item(189,3).Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws.Cells.Item(189,4).Value = "0x00,0x38"
$ws.Cells.Item(189,5).Value = "0x3"
$ws.Cells.Item(189,6).Value = 110
$ws.Cells.Item(189,7).Value = [double]"9.85046333984776e+23"
$ws.Cells.Item(189,8).Value = 69
$ws.Cells.Item(189,9).Value = 3
$ws.Cells.Item(189,1).NumberFormat = $ws.Cells.Item(181,1).NumberFormat()
